$d = $word.ActiveDocument

# --- Change 1: append extra text to the first paragraph ---------------
$para1 = $d.Paragraphs(1)
$r = $para1.Range
$r.Collapse(0)  # wdCollapseEnd -> collapse to end of paragraph text (before the pilcrow)
$r.MoveEnd(1, -1) | Out-Null

$end1 = $para1.Range
$end1.SetRange($end1.End - 1, $end1.End - 1)
$end1.InsertAfter("  ")

$end2 = $d.Range($end1.End, $end1.End)
$end2.InsertAfter([char]40 + "This is a change " + [char]0x2013 + " Ve")
$end2.Font.Color = 12583104

$end3 = $d.Range($end2.End, $end2.End)
$end3.InsertAfter("rsion for branch alternate")
$end3.Font.Color = 12583104

$end4 = $d.Range($end3.End, $end3.End)
$end4.InsertAfter(")")
$end4.Font.Color = 12583104

Write-Output "done change1"
